$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates --------------------------------------------------
# Row 5: helper-plate note gets a pregrown-helper caveat for assay IIc.
$ws.Range("A5").Value = "helper plate (384 w- transparent). For Iic use pregrown helper"
# Row 11: overnight-strain note now excludes IIc.
$ws.Range("A11").Value = "12 Col Overnight Strains (not for IIc)"
# Row 14: the "For IIb" caveat now also excludes IIc.
$ws.Range("A14").Value = "For IIb (and not for IIa, IIc) "
# Row 15: replaced by the new earlier-distribution-to-helper instruction.
$ws.Range("A15").Value = " Add second 384 transp, plate to cart 2, site 1 as next helper"

# --- Formatting updates ------------------------------------------------------
# Rows 2-13: column A switches to Text number format + wrap text (keeps the
# existing thin border / default font). Build the format once on A2, then
# fan it out with a format-only paste so we don't leave extra transient
# cell-style records behind.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").WrapText = $true
$ws.Range("A2").Copy()
$ws.Range("A3:A13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 14 ("For IIb ...") also becomes Text format + wrap text, keeping its
# existing bold font and left/right border.
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").WrapText = $true

# --- Column width -------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 29.76

# --- Selection ------------------------------------------------------------------
$ws.Range("L13").Select()
